$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header in column D
$ws.Range("D1").Value = "is_active"

# Fill column D (rows 2-9) with value 1 (active) for each currency record
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 4).Value = 1
}

# Update selection to match target state (F6)
$ws.Range("F6").Select()
